$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "62.683.83"
$ws.Cells.Item(2,5).Value = "  -0.35%  "
$ws.Cells.Item(3,4).Value = "3.443.53"
$ws.Cells.Item(3,5).Value = "  -0.97%  "
$ws.Cells.Item(4,5).Value = "  -0.04%  "
$ws.Cells.Item(5,4).Value = "'579.63"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = "  -0.54%  "
$ws.Cells.Item(6,4).Value = "'147.76"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = "  +0.48%  "
$ws.Cells.Item(7,5).Value = "  +0.01%  "
$ws.Cells.Item(8,5).Value = "  +0.52%  "
$ws.Cells.Item(9,5).Value = "  +4.15%  "
$ws.Cells.Item(10,5).Value = "  -1.67%  "
$ws.Cells.Item(11,5).Value = "  +2.86%  "
$ws.Cells.Item(12,4).Value = "4.031.32"
$ws.Cells.Item(13,5).Value = "  +2.14%  "
$ws.Cells.Item(14,4).Value = "'28.22"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = "  -5.23%  "
$ws.Cells.Item(15,4).Value = "3.446.36"
$ws.Cells.Item(15,5).Value = "  -1.25%  "
$ws.Cells.Item(16,5).Value = "  -0.30%  "
$ws.Cells.Item(17,4).Value = "62.713.22"
$ws.Cells.Item(17,5).Value = "  -0.31%  "
$ws.Cells.Item(18,4).Value = "'6.37"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value = "  +0.85%  "
$ws.Cells.Item(19,4).Value = "'14.66"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value = "  +2.16%  "
$ws.Cells.Item(20,4).Value = "'9.06"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = "  -2.66%  "
$ws.Cells.Item(21,4).Value = "'386.86"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = "  -0.29%  "
$ws.Cells.Item(22,4).Value = "'75.28"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = "  +0.47%  "
$ws.Cells.Item(23,4).Value = "'0.561"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value = "  -0.45%  "
$ws.Cells.Item(24,5).Value = "  +0.29%  "
$ws.Cells.Item(25,5).Value = "  -1.34%  "
$ws.Cells.Item(26,5).Value = "  -1.71%  "
$ws.Cells.Item(27,5).Value = "  +0.42%  "
$ws.Cells.Item(28,4).Value = "'7.63"
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value = "  -0.56%  "
$ws.Cells.Item(29,4).Value = "'0.999"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value = "  -0.09%  "
$ws.Cells.Item(30,5).Value = "  -3.11%  "
$ws.Cells.Item(31,5).Value = "  -1.17%  "
$ws.Cells.Item(32,5).Value = "  -0.01%  "
$ws.Cells.Item(33,4).Value = "'1.35"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value = "  -4.40%  "
$ws.Cells.Item(34,4).Value = "'23.21"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value = "  -2.35%  "
$ws.Cells.Item(35,4).Value = "'1.63"
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).Value = "  +3.36%  "
$ws.Cells.Item(36,4).Value = "'5.34"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Value = "  +0.58%  "
$ws.Cells.Item(37,4).Value = "'31.86"
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value = "  -0.11%  "
$ws.Cells.Item(38,4).Value = "'6.97"
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value = "  -1.85%  "
$ws.Cells.Item(39,4).Value = "'169.48"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value = "  -1.06%  "
$ws.Cells.Item(40,4).Value = "3.477.68"
$ws.Cells.Item(40,5).Value = "  -1.17%  "
$ws.Cells.Item(41,4).Value = "'0.0773"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = "  +0.51%  "
$ws.Cells.Item(42,5).Value = "  -2.48%  "
$ws.Cells.Item(43,4).Value = "'42.61"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value = "  +0.85%  "
$ws.Cells.Item(44,5).Value = "  -1.05%  "
$ws.Cells.Item(45,4).Value = "'4.36"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = "  -2.46%  "
$ws.Cells.Item(46,5).Value = "  -1.50%  "
$ws.Cells.Item(47,4).Value = "2.567.43"
$ws.Cells.Item(47,5).Value = "  -2.01%  "
$ws.Cells.Item(48,4).Value = "'6.94"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value = "  +2.26%  "
$ws.Cells.Item(49,5).Value = "  -0.79%  "
$ws.Cells.Item(50,4).Value = "'22.59"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = "  -3.68%  "
$ws.Cells.Item(51,4).Value = "'0.999"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = "  -0.15%  "
